$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 46/47: Coin/Link order swap (PaxDollar <-> Decentraland) with new Price/Volume ---
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4638"
$ws.Range("E46").Value = "  -1.00%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9995"
$ws.Range("E47").Value = "  +0.00%  "

# --- Price (D) / Volume(1h) (E) updates for remaining rows ---
# D column values are set as Text (NumberFormat "@") to preserve literal
# formatting (e.g. thousand-dot separators, trailing zeros, leading zeros)
# exactly as scraped, matching the original inline-string cells.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.840.82"
$ws.Range("E2").Value = "  -1.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.857.67"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.16"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5058"
$ws.Range("E7").Value = "  -1.42%  "
$ws.Range("E8").Value = "  -3.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07159"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8903"
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.66"
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.859.44"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07418"
$ws.Range("E13").Value = "  -2.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.35"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.223"
$ws.Range("E15").Value = "  -2.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008503"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.04"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.884.09"
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.014"
$ws.Range("E21").Value = "  -1.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.094.49"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("E23").Value = "  -3.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.429"
$ws.Range("E24").Value = "  -1.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.12"
$ws.Range("E25").Value = "  -2.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.796"
$ws.Range("E26").Value = "  -2.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.81"
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.055"
$ws.Range("E28").Value = "  -3.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.99"
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.630"
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.660"
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09220"
$ws.Range("E32").Value = "  +2.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05083"
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.983"
$ws.Range("E34").Value = "  -3.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7408"
$ws.Range("E35").Value = "  -2.03%  "
$ws.Range("E36").Value = "  -2.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.240"
$ws.Range("E37").Value = "  +6.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.501"
$ws.Range("E38").Value = "  -1.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01988"
$ws.Range("E39").Value = "  -2.74%  "
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5326"
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "119.43"
$ws.Range("E42").Value = "  +4.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.463"
$ws.Range("E43").Value = "  -2.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.363"
$ws.Range("E44").Value = "  -2.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1456"
$ws.Range("E45").Value = "  -1.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.989"
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.558"
$ws.Range("E49").Value = "  -1.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.77"
$ws.Range("E50").Value = "  +0.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.84"
$ws.Range("E51").Value = "  -3.67%  "
